# Add support for ex3400 switches: new "model" column, plus a new device row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for the ex3400 device (junos) right after the ptsw2-floor7 row.
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "10.9.106.11"
$ws.Range("B5").Value = "junos"
$ws.Range("C5").Value = "oren-flr1sw-B1"

# New "model" column (D), with header + values for every data row.
$ws.Range("D1").Value = "model"
$ws.Range("D2").Value = "ex3300"
$ws.Range("D3").Value = "ex3300"
$ws.Range("D4").Value = "ex3300"
$ws.Range("D5").Value = "ex3400"
$ws.Range("D6").Value = "ex3300"
$ws.Range("D7").Value = "ex3400"
$ws.Range("D8").Value = "vsrx"

$ws.Range("B6").Select()
